# Applies the "Added periodic & upfront related scenarios" edit:
#  - ProductLoanInput!B17 changes from "Mifos style" to
#    "Penalties, Fees, Interest, Principal order", with a new left/top
#    aligned style.
#  - ProductLoanInput becomes the active sheet/tab, with B17 selected.

$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")

$cell = $wsInput.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

$wsInput.Activate()
$cell.Select()
